$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 596: FTAI
$ws.Cells.Item(596, 1).Value = "FTAI"
$ws.Cells.Item(596, 2).NumberFormat = "@"
$ws.Cells.Item(596, 2).Value = "2025-08-08"
$ws.Cells.Item(596, 2).ClearFormats()
$ws.Cells.Item(596, 3).Value = 136.06
$ws.Cells.Item(596, 4).Value = 276.01
$ws.Cells.Item(596, 5).Value = 102.86
$ws.Cells.Item(596, 6).Value = 11.21528

# Row 602: MOG-B
$ws.Cells.Item(602, 1).Value = "MOG-B"
$ws.Cells.Item(602, 2).NumberFormat = "@"
$ws.Cells.Item(602, 2).Value = "2025-08-08"
$ws.Cells.Item(602, 2).ClearFormats()
$ws.Cells.Item(602, 3).Value = 194.99
$ws.Cells.Item(602, 4).Value = 295.0
$ws.Cells.Item(602, 5).Value = 51.29
$ws.Cells.Item(602, 6).Value = 6.99596

# Row 603: AKAM
$ws.Cells.Item(603, 1).Value = "AKAM"
$ws.Cells.Item(603, 2).NumberFormat = "@"
$ws.Cells.Item(603, 2).Value = "2025-08-08"
$ws.Cells.Item(603, 2).ClearFormats()
$ws.Cells.Item(603, 3).Value = 70.53
$ws.Cells.Item(603, 4).Value = 99.35
$ws.Cells.Item(603, 5).Value = 40.86
$ws.Cells.Item(603, 6).Value = 3.37612

# Row 604: KRMN
$ws.Cells.Item(604, 1).Value = "KRMN"
$ws.Cells.Item(604, 2).NumberFormat = "@"
$ws.Cells.Item(604, 2).Value = "2025-08-08"
$ws.Cells.Item(604, 2).ClearFormats()
$ws.Cells.Item(604, 3).Value = 45.78
$ws.Cells.Item(604, 4).Value = 110.93
$ws.Cells.Item(604, 5).Value = 142.31
$ws.Cells.Item(604, 6).Value = 4.87571

# Row 605: ADI
$ws.Cells.Item(605, 1).Value = "ADI"
$ws.Cells.Item(605, 2).NumberFormat = "@"
$ws.Cells.Item(605, 2).Value = "2025-08-08"
$ws.Cells.Item(605, 2).ClearFormats()
$ws.Cells.Item(605, 3).Value = 223.95
$ws.Cells.Item(605, 4).Value = 318.7
$ws.Cells.Item(605, 5).Value = 42.31
$ws.Cells.Item(605, 6).Value = 6.15998

# Row 606: BIO-B
$ws.Cells.Item(606, 1).Value = "BIO-B"
$ws.Cells.Item(606, 2).NumberFormat = "@"
$ws.Cells.Item(606, 2).Value = "2025-08-08"
$ws.Cells.Item(606, 2).ClearFormats()
$ws.Cells.Item(606, 3).Value = 247.05
$ws.Cells.Item(606, 4).Value = 312.75
$ws.Cells.Item(606, 5).Value = 26.59
$ws.Cells.Item(606, 6).Value = 0.4732

# Row 607: TEVA
$ws.Cells.Item(607, 1).Value = "TEVA"
$ws.Cells.Item(607, 2).NumberFormat = "@"
$ws.Cells.Item(607, 2).Value = "2025-08-08"
$ws.Cells.Item(607, 2).ClearFormats()
$ws.Cells.Item(607, 3).Value = 16.37
$ws.Cells.Item(607, 4).Value = 33.07
$ws.Cells.Item(607, 5).Value = 102.02
$ws.Cells.Item(607, 6).Value = 0.43058

# Row 608: COR
$ws.Cells.Item(608, 1).Value = "COR"
$ws.Cells.Item(608, 2).NumberFormat = "@"
$ws.Cells.Item(608, 2).Value = "2025-08-08"
$ws.Cells.Item(608, 2).ClearFormats()
$ws.Cells.Item(608, 3).Value = 285.13
$ws.Cells.Item(608, 4).Value = 352.1
$ws.Cells.Item(608, 5).Value = 23.49
$ws.Cells.Item(608, 6).Value = 2.00138

# Row 609: BNS
$ws.Cells.Item(609, 1).Value = "BNS"
$ws.Cells.Item(609, 2).NumberFormat = "@"
$ws.Cells.Item(609, 2).Value = "2025-08-08"
$ws.Cells.Item(609, 2).ClearFormats()
$ws.Cells.Item(609, 3).Value = 56.0
$ws.Cells.Item(609, 4).Value = 76.42
$ws.Cells.Item(609, 5).Value = 36.46
$ws.Cells.Item(609, 6).Value = 0.83523

# Row 610: TIGO
$ws.Cells.Item(610, 1).Value = "TIGO"
$ws.Cells.Item(610, 2).NumberFormat = "@"
$ws.Cells.Item(610, 2).Value = "2025-08-08"
$ws.Cells.Item(610, 2).ClearFormats()
$ws.Cells.Item(610, 3).Value = 42.05
$ws.Cells.Item(610, 4).Value = 62.22
$ws.Cells.Item(610, 5).Value = 47.97
$ws.Cells.Item(610, 6).Value = 1.76795

# Row 611: MLI
$ws.Cells.Item(611, 1).Value = "MLI"
$ws.Cells.Item(611, 2).NumberFormat = "@"
$ws.Cells.Item(611, 2).Value = "2025-08-08"
$ws.Cells.Item(611, 2).ClearFormats()
$ws.Cells.Item(611, 3).Value = 89.08
$ws.Cells.Item(611, 4).Value = 137.38
$ws.Cells.Item(611, 5).Value = 54.22
$ws.Cells.Item(611, 6).Value = 3.23044

# Row 612: SAN
$ws.Cells.Item(612, 1).Value = "SAN"
$ws.Cells.Item(612, 2).NumberFormat = "@"
$ws.Cells.Item(612, 2).Value = "2025-08-08"
$ws.Cells.Item(612, 2).ClearFormats()
$ws.Cells.Item(612, 3).Value = 9.31
$ws.Cells.Item(612, 4).Value = 12.81
$ws.Cells.Item(612, 5).Value = 37.59
$ws.Cells.Item(612, 6).Value = 0.26373

# Row 613: SLF
$ws.Cells.Item(613, 1).Value = "SLF"
$ws.Cells.Item(613, 2).NumberFormat = "@"
$ws.Cells.Item(613, 2).Value = "2025-08-08"
$ws.Cells.Item(613, 2).ClearFormats()
$ws.Cells.Item(613, 3).Value = 56.54
$ws.Cells.Item(613, 4).Value = 63.79
$ws.Cells.Item(613, 5).Value = 12.82
$ws.Cells.Item(613, 6).Value = 0.28712

# Row 614: GH
$ws.Cells.Item(614, 1).Value = "GH"
$ws.Cells.Item(614, 2).NumberFormat = "@"
$ws.Cells.Item(614, 2).Value = "2025-08-08"
$ws.Cells.Item(614, 2).ClearFormats()
$ws.Cells.Item(614, 3).Value = 54.23
$ws.Cells.Item(614, 4).Value = 110.48
$ws.Cells.Item(614, 5).Value = 103.72
$ws.Cells.Item(614, 6).Value = 0.89182

# Row 615: BBIO
$ws.Cells.Item(615, 1).Value = "BBIO"
$ws.Cells.Item(615, 2).NumberFormat = "@"
$ws.Cells.Item(615, 2).Value = "2025-08-08"
$ws.Cells.Item(615, 2).ClearFormats()
$ws.Cells.Item(615, 3).Value = 46.58
$ws.Cells.Item(615, 4).Value = 79.1
$ws.Cells.Item(615, 5).Value = 69.82
$ws.Cells.Item(615, 6).Value = 0.51568

# Row 616: BE
$ws.Cells.Item(616, 1).Value = "BE"
$ws.Cells.Item(616, 2).NumberFormat = "@"
$ws.Cells.Item(616, 2).Value = "2025-08-08"
$ws.Cells.Item(616, 2).ClearFormats()
$ws.Cells.Item(616, 3).Value = 36.8
$ws.Cells.Item(616, 4).Value = 156.51
$ws.Cells.Item(616, 5).Value = 325.3
$ws.Cells.Item(616, 6).Value = 8.80753

# Row 617: ARWR
$ws.Cells.Item(617, 1).Value = "ARWR"
$ws.Cells.Item(617, 2).NumberFormat = "@"
$ws.Cells.Item(617, 2).Value = "2025-08-08"
$ws.Cells.Item(617, 2).ClearFormats()
$ws.Cells.Item(617, 3).Value = 16.42
$ws.Cells.Item(617, 4).Value = 69.12
$ws.Cells.Item(617, 5).Value = 320.95
$ws.Cells.Item(617, 6).Value = 0.34305

# Row 618: EDU
$ws.Cells.Item(618, 1).Value = "EDU"
$ws.Cells.Item(618, 2).NumberFormat = "@"
$ws.Cells.Item(618, 2).Value = "2025-08-08"
$ws.Cells.Item(618, 2).ClearFormats()
$ws.Cells.Item(618, 3).Value = 45.58
$ws.Cells.Item(618, 4).Value = 59.46
$ws.Cells.Item(618, 5).Value = 30.45
$ws.Cells.Item(618, 6).Value = 0.47873

# Row 619: INTC
$ws.Cells.Item(619, 1).Value = "INTC"
$ws.Cells.Item(619, 2).NumberFormat = "@"
$ws.Cells.Item(619, 2).Value = "2025-08-08"
$ws.Cells.Item(619, 2).ClearFormats()
$ws.Cells.Item(619, 3).Value = 19.95
$ws.Cells.Item(619, 4).Value = 48.66
$ws.Cells.Item(619, 5).Value = 143.91
$ws.Cells.Item(619, 6).Value = 0.39579

# Row 620: FNV
$ws.Cells.Item(620, 1).Value = "FNV"
$ws.Cells.Item(620, 2).NumberFormat = "@"
$ws.Cells.Item(620, 2).Value = "2025-08-08"
$ws.Cells.Item(620, 2).ClearFormats()
$ws.Cells.Item(620, 3).Value = 171.59
$ws.Cells.Item(620, 4).Value = 261.77
$ws.Cells.Item(620, 5).Value = 52.56
$ws.Cells.Item(620, 6).Value = 9.51144

# Row 621: VIK
$ws.Cells.Item(621, 1).Value = "VIK"
$ws.Cells.Item(621, 2).NumberFormat = "@"
$ws.Cells.Item(621, 2).Value = "2025-08-08"
$ws.Cells.Item(621, 2).ClearFormats()
$ws.Cells.Item(621, 3).Value = 56.58
$ws.Cells.Item(621, 4).Value = 74.29
$ws.Cells.Item(621, 5).Value = 31.3
$ws.Cells.Item(621, 6).Value = 0.12575

# Row 622: SPG
$ws.Cells.Item(622, 1).Value = "SPG"
$ws.Cells.Item(622, 2).NumberFormat = "@"
$ws.Cells.Item(622, 2).Value = "2025-08-08"
$ws.Cells.Item(622, 2).ClearFormats()
$ws.Cells.Item(622, 3).Value = 164.93
$ws.Cells.Item(622, 4).Value = 189.8
$ws.Cells.Item(622, 5).Value = 15.08
$ws.Cells.Item(622, 6).Value = 0.43519
